$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '42.921.07'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = "'" + '2.552.14'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'" + '303.68'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").Value = "'" + '98.25'
$ws.Range("E6").Value = '  +4.12%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'" + '0.545'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = "'" + '36.63'
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").Value = "'" + '0.0818'
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = "'" + '0.116'
$ws.Range("E12").Value = '  +6.00%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'" + '7.62'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").Value = "'" + '2.941.92'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = "'" + '2.592.65'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = "'" + '14.85'
$ws.Range("E16").Value = '  +4.77%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = "'" + '0.879'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").Value = "'" + '43.135.62'
$ws.Range("D19").Value = "'" + '13.65'
$ws.Range("E19").Value = '  +4.95%  '
$ws.Range("D20").Value = "'" + '0.0₃0985'
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").Value = "'" + '6.61'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = "'" + '71.90'
$ws.Range("D23").Value = "'" + '254.59'
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("D24").Value = "'" + '2.96'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").Value = "'" + '27.90'
$ws.Range("E26").Value = '  -6.02%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = "'" + '10.13'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = "'" + '37.81'
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").Value = "'" + '159.16'
$ws.Range("E32").Value = '  +3.04%  '
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("D35").Value = "'" + '0.0804'
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("D36").Value = "'" + '3.32'
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("D37").Value = "'" + '18.86'
$ws.Range("E37").Value = '  +12.76%  '
$ws.Range("D38").Value = "'" + '25.78'
$ws.Range("E38").Value = '  +10.70%  '
$ws.Range("E39").Value = '  -1.40%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("B41").Value = 'ApeXProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D41").Value = "'" + '2.10'
$ws.Range("E41").Value = '  +32.25%  '
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").Value = "'" + '3.43'
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = "'" + '2.097.97'
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").Value = "'" + '86.35'
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("E48").Value = '  +3.08%  '
$ws.Range("D49").Value = "'" + '2.799.52'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = "'" + '74.71'
$ws.Range("E50").Value = '  +7.53%  '
$ws.Range("D51").Value = "'" + '103.60'
$ws.Range("E51").Value = '  -0.89%  '
